$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = 0.4450587250587252
$ws.Range("F4").Value = 0.9815481935481936
$ws.Range("H4").Value = 0.398956550956551
$ws.Range("J4").Value = 0.2904358980202829
$ws.Range("K4").Value = -0.03217824817824818
$ws.Range("N4").Value = -0.7553534393534395
$ws.Range("P4").Value = -0.5370557490557492
$ws.Range("C5").Value = -0.02522546522546523
$ws.Range("F5").Value = 0.07346710946710948
$ws.Range("H5").Value = 0.1268329868329869
$ws.Range("J5").Value = 0.1846338707267886
$ws.Range("K5").Value = 0.007497619497619499
$ws.Range("N5").Value = -0.01155748755748756
$ws.Range("P5").Value = -0.1731345411345411
$ws.Range("C6").Value = -0.01757182157182157
$ws.Range("F6").Value = 0.1024356664356665
$ws.Range("H6").Value = 0.05476237876237876
$ws.Range("J6").Value = 0.01137926975478567
$ws.Range("K6").Value = -0.1052066492066492
$ws.Range("N6").Value = 0.3132213252213252
$ws.Range("P6").Value = 0.1284032004032004
$ws.Range("C7").Value = 0.07966122766122767
$ws.Range("F7").Value = -0.1074869274869275
$ws.Range("H7").Value = -0.08355268755268756
$ws.Range("J7").Value = -0.0192572118926999
$ws.Range("K7").Value = 0.1335063375063375
$ws.Range("N7").Value = -0.4482812562812563
$ws.Range("P7").Value = -0.1847998247998248
$ws.Range("C8").Value = 0.1578608298608299
$ws.Range("F8").Value = 0.124034128034128
$ws.Range("H8").Value = 0.907915327915328
$ws.Range("J8").Value = -0.3274756882823862
$ws.Range("K8").Value = -0.02194046194046194
$ws.Range("N8").Value = 0.05802094602094603
$ws.Range("P8").Value = 0.4858709218709218
$ws.Range("C9").Value = 0.01345491745491745
$ws.Range("F9").Value = 0.04774168774168774
$ws.Range("H9").Value = 0.03732184932184932
$ws.Range("J9").Value = 0.02300513976430939
$ws.Range("K9").Value = -0.02924674124674125
$ws.Range("N9").Value = 0.07657105657105658
$ws.Range("P9").Value = -0.02185917385917386
$ws.Range("C10").Value = -0.01268656868656869
$ws.Range("F10").Value = 0.04083012483012483
$ws.Range("H10").Value = -0.04444772044772045
$ws.Range("J10").Value = 0.04542375492273887
$ws.Range("K10").Value = -0.02535238935238936
$ws.Range("N10").Value = -0.006716406716406716
$ws.Range("P10").Value = -0.01201078801078801
$ws.Range("C11").Value = -0.0231927111927112
$ws.Range("F11").Value = -0.04084812484812485
$ws.Range("H11").Value = -0.03075438675438676
$ws.Range("J11").Value = 0.001185511224633096
$ws.Range("K11").Value = 0.0704925224925225
$ws.Range("N11").Value = -0.0007574887574887575
$ws.Range("P11").Value = -0.0290938250938251
$ws.Range("C12").Value = 0.01755983355983356
$ws.Range("F12").Value = -0.003258267258267258
$ws.Range("H12").Value = -0.03247745647745648
$ws.Range("J12").Value = -0.01952971417419476
$ws.Range("K12").Value = 0.0462946182946183
$ws.Range("N12").Value = -0.05350919350919351
$ws.Range("P12").Value = -0.08013543213543213
$ws.Range("C13").Value = 0.03985986385986386
$ws.Range("F13").Value = 0.03533617133617133
$ws.Range("H13").Value = 0.02935520935520936
$ws.Range("J13").Value = 0.0309902840129634
$ws.Range("K13").Value = -0.0110995430995431
$ws.Range("N13").Value = -0.0474975594975595
$ws.Range("P13").Value = -0.007846975846975848
$ws.Range("C14").Value = 0.4944883224883225
$ws.Range("F14").Value = -0.04154450954450954
$ws.Range("H14").Value = -0.0105029145029145
$ws.Range("J14").Value = 0.0202871069565815
$ws.Range("K14").Value = -0.622924342924343
$ws.Range("N14").Value = 0.08232916632916634
$ws.Range("P14").Value = 0.07736422136422137
$ws.Range("C15").Value = 0.5499836499836501
$ws.Range("F15").Value = -0.03117846717846719
$ws.Range("H15").Value = -0.06312644712644713
$ws.Range("J15").Value = 0.02399873479069306
$ws.Range("K15").Value = 0.7433979953979954
$ws.Range("N15").Value = -0.01178776778776779
$ws.Range("P15").Value = -0.04887460887460888
$ws.Range("C16").Value = -0.04426584826584827
$ws.Range("F16").Value = -0.01595469995469995
$ws.Range("H16").Value = -0.006861606861606862
$ws.Range("J16").Value = -0.01002586035671376
$ws.Range("K16").Value = 0.02162897762897763
$ws.Range("N16").Value = 0.02258873858873859
$ws.Range("P16").Value = 0.01402796602796603
$ws.Range("C17").Value = 0.03324890124890125
$ws.Range("F17").Value = 0.002408582408582409
$ws.Range("H17").Value = 0.01225959625959626
$ws.Range("J17").Value = -0.03404621916974443
$ws.Range("K17").Value = 0.02187738987738988
$ws.Range("N17").Value = -0.01244346044346044
$ws.Range("P17").Value = -0.0432970512970513
$ws.Range("C18").Value = 0.081993573993574
$ws.Range("F18").Value = 0.04557774957774958
$ws.Range("H18").Value = 0.0008977328977328978
$ws.Range("J18").Value = 0.02816583509530767
$ws.Range("K18").Value = -0.0004415404415404417
$ws.Range("N18").Value = -0.0536956016956017
$ws.Range("P18").Value = -0.01426904626904627
$ws.Range("C19").Value = -0.2105943905943906
$ws.Range("F19").Value = 0.00868943668943669
$ws.Range("H19").Value = -0.003231411231411232
$ws.Range("J19").Value = 0.02864605759137751
$ws.Range("K19").Value = -0.09644532044532045
$ws.Range("N19").Value = -0.01376244176244176
$ws.Range("P19").Value = -0.02062486462486463
$ws.Range("C20").Value = 0.04410844410844411
$ws.Range("F20").Value = 0.001672585672585673
$ws.Range("H20").Value = 0.01854306654306654
$ws.Range("J20").Value = 0.0103137586541127
$ws.Range("K20").Value = 0.02843884043884044
$ws.Range("N20").Value = 0.009667641667641669
$ws.Range("P20").Value = 0.005071241071241072
$ws.Range("C21").Value = -0.004607380607380607
$ws.Range("F21").Value = 0.01538705138705139
$ws.Range("H21").Value = -0.001155661155661156
$ws.Range("J21").Value = 0.02578636463731469
$ws.Range("K21").Value = -0.004811512811512811
$ws.Range("N21").Value = -0.01106756306756307
$ws.Range("P21").Value = -0.0157917757917758
$ws.Range("C22").Value = -0.03126273126273126
$ws.Range("F22").Value = -0.006104226104226105
$ws.Range("H22").Value = 0.0183978423978424
$ws.Range("J22").Value = -0.04032669565747664
$ws.Range("K22").Value = 0.005055377055377056
$ws.Range("N22").Value = 0.01157055557055557
$ws.Range("P22").Value = 0.004980640980640981
$ws.Range("C23").Value = -0.01001806601806602
$ws.Range("F23").Value = -0.01162511962511963
$ws.Range("H23").Value = 0.008321252321252321
$ws.Range("J23").Value = 0.04679035833444019
$ws.Range("K23").Value = 0.003523863523863525
$ws.Range("N23").Value = 0.06616223416223416
$ws.Range("P23").Value = 0.04403445203445204
$ws.Range("C24").Value = 0.004313596313596313
$ws.Range("F24").Value = 0.01884685884685885
$ws.Range("H24").Value = -0.03162809562809563
$ws.Range("J24").Value = 0.01761580819712988
$ws.Range("K24").Value = 0.04577704577704578
$ws.Range("N24").Value = -0.04873269673269674
$ws.Range("P24").Value = -0.04377267177267178
$ws.Range("C25").Value = -0.001007653007653008
$ws.Range("F25").Value = -0.03705560505560506
$ws.Range("H25").Value = -0.02398396798396799
$ws.Range("J25").Value = -0.02557333441725447
$ws.Range("K25").Value = 0.005005853005853007
$ws.Range("N25").Value = -0.022003042003042
$ws.Range("P25").Value = 0.02104720504720505
$ws.Range("C26").Value = 0.04556800556800557
$ws.Range("F26").Value = 0.03614800814800815
$ws.Range("H26").Value = 0.0171989811989812
$ws.Range("J26").Value = 0.02191193863503262
$ws.Range("K26").Value = -0.08105319305319306
$ws.Range("N26").Value = 0.0002172362172362172
$ws.Range("P26").Value = -0.003624255624255624
$ws.Range("C27").Value = 0.0004721644721644722
$ws.Range("F27").Value = 0.02612502212502213
$ws.Range("H27").Value = -0.0040991800991801
$ws.Range("J27").Value = -0.03363122874105931
$ws.Range("K27").Value = -0.02553369753369754
$ws.Range("N27").Value = -0.01381768981768982
$ws.Range("P27").Value = -0.009276213276213276
$ws.Range("C28").Value = -0.005404805404805406
$ws.Range("F28").Value = -0.001617001617001617
$ws.Range("H28").Value = -0.02635341835341836
$ws.Range("J28").Value = -0.0003738723862101752
$ws.Range("K28").Value = -0.005721377721377722
$ws.Range("N28").Value = -0.01993636393636394
$ws.Range("P28").Value = -0.002308010308010308
$ws.Range("C29").Value = -0.02977214977214978
$ws.Range("F29").Value = -0.05065675465675466
$ws.Range("H29").Value = -0.02243819843819844
$ws.Range("J29").Value = -0.03631537351378087
$ws.Range("K29").Value = -0.009992697992697993
$ws.Range("N29").Value = 0.02819384819384819
$ws.Range("P29").Value = -0.01042933042933043
$ws.Range("C30").Value = 0.01090565890565891
$ws.Range("F30").Value = 0.006529998529998531
$ws.Range("H30").Value = 0.02907527307527308
$ws.Range("J30").Value = 0.01991552057273276
$ws.Range("K30").Value = 0.02343598743598744
$ws.Range("N30").Value = -0.02546445746445746
$ws.Range("P30").Value = 0.04133420933420934
